# Weekly refresh of the "Hortaliza, Terminal Hortofrutícola Agro Chillán -
# Zanahoria" series: a new week's observation is inserted at the top of the
# data block (row 72) and every older observation shifts down one row, with
# the previously-last row (192) now landing on the new row 193.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 72:192 down to 73:193, leaving a blank row 72 behind
# (format of the date column is carried along automatically by Excel).
$ws.Rows("72:72").Insert()

# Populate the freshly inserted row with this week's record.
$ws.Range("A72").Value = 7
$ws.Range("B72").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C72").Value = "Ñuble"
$ws.Range("D72").Value = 44495
$ws.Range("E72").Value = 16
$ws.Range("F72").Value = 100114013
$ws.Range("G72").Value = "Zanahoria"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 100
$ws.Range("K72").Value = 8500
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = 8750
$ws.Range("N72").Value = "$/saco 20 kilos"
$ws.Range("O72").Value = "Región de Ñuble"
$ws.Range("P72").Value = 438
$ws.Range("Q72").Value = 20
$ws.Range("R72").Value = "Hortaliza"
